# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 212.2
$ws.Range("I2").Value = 140.25
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 140.25
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -27.25
$ws.Range("N2").Value = -726
$ws.Range("H18").Value = 1234.6666
$ws.Range("I18").Value = 1173.7142
$ws.Range("J18").Value = 1448
$ws.Range("K18").Value = 1173.7142
$ws.Range("L18").Value = 1448
$ws.Range("M18").Value = -889.7141999999999
$ws.Range("N18").Value = -2016
$ws.Range("H41").Value = 1011.5833
$ws.Range("I41").Value = 830.8182
$ws.Range("K41").Value = 830.8182
$ws.Range("M41").Value = -390.8182
$ws.Range("H69").Value = 7760.2173
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7760.2173
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 23280.6519
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -25028.6519
$ws.Range("H70").Value = 2941
$ws.Range("I70").Value = 1803.6666
$ws.Range("J70").Value = 3580.75
$ws.Range("K70").Value = 5410.9998
$ws.Range("L70").Value = 10742.25
$ws.Range("M70").Value = -5140.9998
$ws.Range("N70").Value = -11282.25
$ws.Range("H72").Value = 7760.2173
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7760.2173
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 69841.95570000001
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -78577.95570000001
$ws.Range("H73").Value = 2941
$ws.Range("I73").Value = 1803.6666
$ws.Range("J73").Value = 3580.75
$ws.Range("K73").Value = 5410.9998
$ws.Range("L73").Value = 10742.25
$ws.Range("M73").Value = -4474.9998
$ws.Range("N73").Value = -12614.25
$ws.Range("H86").Value = 2151.2
$ws.Range("I86").Value = 1928.25
$ws.Range("J86").Value = 2299.8333
$ws.Range("K86").Value = 1928.25
$ws.Range("L86").Value = 2299.8333
$ws.Range("M86").Value = -805.25
$ws.Range("N86").Value = -4545.8333
$ws.Range("H89").Value = 2151.2
$ws.Range("I89").Value = 1928.25
$ws.Range("J89").Value = 2299.8333
$ws.Range("K89").Value = 9641.25
$ws.Range("L89").Value = 11499.1665
$ws.Range("M89").Value = -4025.25
$ws.Range("N89").Value = -22731.1665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3047.8948
$ws.Range("I86").Value = 1391.2
$ws.Range("J86").Value = 4888.6665
$ws.Range("K86").Value = 1391.2
$ws.Range("L86").Value = 4888.6665
$ws.Range("M86").Value = -268.2
$ws.Range("N86").Value = -7134.6665
$ws.Range("H89").Value = 3047.8948
$ws.Range("I89").Value = 1391.2
$ws.Range("J89").Value = 4888.6665
$ws.Range("K89").Value = 6956
$ws.Range("L89").Value = 24443.3325
$ws.Range("M89").Value = -1340
$ws.Range("N89").Value = -35675.3325
$ws.Range("H105").Value = 9092706
$ws.Range("I105").Value = 9092706
$ws.Range("K105").Value = 9092706
$ws.Range("M105").Value = -9090959

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1771.1177
$ws.Range("I132").Value = 1694.3125
$ws.Range("K132").Value = 5082.9375
$ws.Range("M132").Value = -2552.9375
$ws.Range("H141").Value = 54213
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 54213
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 54213
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -64573

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1269.6666
$ws.Range("I5").Value = 1094.4
$ws.Range("J5").Value = 1394.8572
$ws.Range("K5").Value = 3283.2
$ws.Range("L5").Value = 4184.571599999999
$ws.Range("M5").Value = -3171.2
$ws.Range("N5").Value = -4408.571599999999
$ws.Range("H60").Value = 1240.2354
$ws.Range("I60").Value = 230.625
$ws.Range("K60").Value = 691.875
$ws.Range("M60").Value = -440.875
$ws.Range("H97").Value = 398.5
$ws.Range("I97").Value = 397
$ws.Range("K97").Value = 1191
$ws.Range("M97").Value = -695
$ws.Range("H132").Value = 2432.4285
$ws.Range("I132").Value = 1455.6
$ws.Range("K132").Value = 13100.4
$ws.Range("M132").Value = -10570.4
$ws.Range("H135").Value = 1269.6666
$ws.Range("I135").Value = 1094.4
$ws.Range("J135").Value = 1394.8572
$ws.Range("K135").Value = 9849.6
$ws.Range("L135").Value = 12553.7148
$ws.Range("M135").Value = -7314.6
$ws.Range("N135").Value = -17623.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 37777.4
$ws.Range("J24").Value = 37777.4
$ws.Range("L24").Value = 37777.4
$ws.Range("N24").Value = -38123.4
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("M80").Value = -2
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("M83").Value = -8
$ws.Range("H132").Value = 1574.3334
$ws.Range("I132").Value = 1574.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4723.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2193.0002
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 909
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 1044.3334
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 1044.3334
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -1634.3334
$ws.Range("H27").Value = 909
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 1044.3334
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 1044.3334
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -1258.3334
$ws.Range("H46").Value = 6517.647
$ws.Range("I46").Value = 4257.143
$ws.Range("K46").Value = 4257.143
$ws.Range("M46").Value = -4069.143
$ws.Range("H68").Value = 5946
$ws.Range("I68").Value = 4594.6665
$ws.Range("K68").Value = 4594.6665
$ws.Range("M68").Value = -3845.6665
$ws.Range("H71").Value = 5946
$ws.Range("I71").Value = 4594.6665
$ws.Range("K71").Value = 22973.3325
$ws.Range("M71").Value = -19229.3325
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
